# Attendance workbook update: mark attendance (value 1) for several
# students on the 2019-12-07 week (column T) and a couple of additional
# marks on the 2019-10-12 week (column L). Finally move the active
# selection to M3 (matches the saved sheetView selection in the file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column T (week of 2019-12-07) newly marked present
$ws.Range("T2").Value = 1
$ws.Range("T5").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("T7").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("T12").Value = 1
$ws.Range("T17").Value = 1
$ws.Range("T18").Value = 1
$ws.Range("T19").Value = 1

# Column L (week of 2019-10-12) newly marked present
$ws.Range("L15").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("L20").Value = 1

# Restore the saved selection/active cell recorded in the sheet view
$ws.Range("M3").Select()
